# Auto-generated Excel COM-interop script
# Applies scheduled-runner market price/profit updates to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 9200.166999999999
$ws.Range("J6").Value = 1040.2
$ws.Range("L6").Value = 3120.6
$ws.Range("N6").Value = -3344.6
$ws.Range("H18").Value = 405
$ws.Range("H38").Value = 238.83333
$ws.Range("J38").Value = 900
$ws.Range("L38").Value = 2700
$ws.Range("N38").Value = -3444
$ws.Range("H41").Value = 1812.2106
$ws.Range("I41").Value = 1452
$ws.Range("J41").Value = 3733.3333
$ws.Range("K41").Value = 1452
$ws.Range("L41").Value = 3733.3333
$ws.Range("M41").Value = -1012
$ws.Range("N41").Value = -4613.3333
$ws.Range("H58").Value = 2768.1333
$ws.Range("J58").Value = 5202.4287
$ws.Range("L58").Value = 15607.2861
$ws.Range("N58").Value = -15907.2861
$ws.Range("H112").Value = 1678.7407
$ws.Range("J112").Value = 1708.7693
$ws.Range("L112").Value = 5126.3079
$ws.Range("N112").Value = -7342.3079
$ws.Range("H116").Value = 2957.9473
$ws.Range("J116").Value = 2666.4
$ws.Range("L116").Value = 2666.4
$ws.Range("N116").Value = -9550.4
$ws.Range("H129").Value = 871.75
$ws.Range("J129").Value = 1000
$ws.Range("L129").Value = 3000
$ws.Range("N129").Value = -13000
$ws.Range("H138").Value = 3574.7612
$ws.Range("I138").Value = 2241.9285
$ws.Range("J138").Value = 3926.83
$ws.Range("K138").Value = 6725.7855
$ws.Range("L138").Value = 11780.49
$ws.Range("M138").Value = -1585.7855
$ws.Range("N138").Value = -22060.49

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2746.9285
$ws.Range("I2").Value = 2094.9
$ws.Range("K2").Value = 2094.9
$ws.Range("M2").Value = -1981.9
$ws.Range("H32").Value = 23908.316
$ws.Range("I32").Value = 21721.117
$ws.Range("J32").Value = 42499.5
$ws.Range("K32").Value = 21721.117
$ws.Range("L32").Value = 42499.5
$ws.Range("M32").Value = -21434.117
$ws.Range("N32").Value = -43073.5
$ws.Range("H116").Value = 2746.9285
$ws.Range("I116").Value = 2094.9
$ws.Range("K116").Value = 2094.9
$ws.Range("M116").Value = 199.0999999999999
$ws.Range("H132").Value = 2986.182
$ws.Range("I132").Value = 2055.2104
$ws.Range("J132").Value = 4249.643
$ws.Range("K132").Value = 6165.6312
$ws.Range("L132").Value = 12748.929
$ws.Range("M132").Value = -3635.6312
$ws.Range("N132").Value = -17808.929

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2746.9285
$ws.Range("I3").Value = 2094.9
$ws.Range("K3").Value = 2094.9
$ws.Range("M3").Value = -1980.9
$ws.Range("H134").Value = 1488.2727
$ws.Range("I134").Value = 1084.1177
$ws.Range("K134").Value = 3252.3531
$ws.Range("M134").Value = -717.3531000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 950
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 950
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 950
$ws.Range("N6").Value = -1176
$ws.Range("M6").ClearContents()
$ws.Range("H7").Value = 207.63637
$ws.Range("I7").Value = 97
$ws.Range("J7").Value = 401.25
$ws.Range("K7").Value = 97
$ws.Range("L7").Value = 401.25
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = -627.25
$ws.Range("H31").Value = 1596.2142
$ws.Range("I31").Value = 1343.65
$ws.Range("J31").Value = 3111.6
$ws.Range("K31").Value = 1343.65
$ws.Range("L31").Value = 3111.6
$ws.Range("M31").Value = -1048.65
$ws.Range("N31").Value = -3701.6
$ws.Range("H34").Value = 1596.2142
$ws.Range("I34").Value = 1343.65
$ws.Range("J34").Value = 3111.6
$ws.Range("K34").Value = 1343.65
$ws.Range("L34").Value = 3111.6
$ws.Range("M34").Value = -1141.65
$ws.Range("N34").Value = -3515.6
$ws.Range("H62").Value = 14288071
$ws.Range("I62").Value = 2538.4614
$ws.Range("K62").Value = 2538.4614
$ws.Range("M62").Value = -1914.4614
$ws.Range("H65").Value = 14288071
$ws.Range("I65").Value = 2538.4614
$ws.Range("K65").Value = 12692.307
$ws.Range("M65").Value = -9572.307000000001
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 32000
$ws.Range("J74").Value = 32000
$ws.Range("L74").Value = 32000
$ws.Range("N74").Value = -33748
$ws.Range("H77").Value = 32000
$ws.Range("J77").Value = 32000
$ws.Range("L77").Value = 96000
$ws.Range("N77").Value = -104736
$ws.Range("H86").Value = 2172466.2
$ws.Range("I86").Value = 3046422.8
$ws.Range("J86").Value = 36127.89
$ws.Range("K86").Value = 3046422.8
$ws.Range("L86").Value = 36127.89
$ws.Range("M86").Value = -3045299.8
$ws.Range("N86").Value = -38373.89
$ws.Range("H89").Value = 2172466.2
$ws.Range("I89").Value = 3046422.8
$ws.Range("J89").Value = 36127.89
$ws.Range("K89").Value = 15232114
$ws.Range("L89").Value = 180639.45
$ws.Range("M89").Value = -15226498
$ws.Range("N89").Value = -191871.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1738
$ws.Range("I69").Value = 1365.6
$ws.Range("J69").Value = 2203.5
$ws.Range("K69").Value = 4096.799999999999
$ws.Range("L69").Value = 6610.5
$ws.Range("M69").Value = -3285.799999999999
$ws.Range("N69").Value = -8232.5
$ws.Range("H72").Value = 1738
$ws.Range("I72").Value = 1365.6
$ws.Range("J72").Value = 2203.5
$ws.Range("K72").Value = 12290.4
$ws.Range("L72").Value = 19831.5
$ws.Range("M72").Value = -8234.4
$ws.Range("N72").Value = -27943.5
$ws.Range("H92").Value = 1066.8334
$ws.Range("I92").Value = 1066.8334
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 3200.5002
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1952.5002
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 730.375
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H127").Value = 2315
$ws.Range("J127").Value = 2315
$ws.Range("L127").Value = 6945
$ws.Range("N127").Value = -16865

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10470.5
$ws.Range("I99").Value = 10470.5
$ws.Range("K99").Value = 10470.5
$ws.Range("M99").Value = -8224.5
$ws.Range("H104").Value = 54166.5
$ws.Range("J104").Value = 54166.5
$ws.Range("L104").Value = 54166.5
$ws.Range("N104").Value = -61154.5
$ws.Range("H113").Value = 1289.25
$ws.Range("I113").Value = 1306.4546
$ws.Range("K113").Value = 1306.4546
$ws.Range("M113").Value = 863.5454
$ws.Range("H116").Value = 45371
$ws.Range("J116").Value = 45371
$ws.Range("L116").Value = 45371
$ws.Range("N116").Value = -54549
$ws.Range("H132").Value = 6214.5864
$ws.Range("I132").Value = 7129.5713
$ws.Range("J132").Value = 3812.75
$ws.Range("K132").Value = 21388.7139
$ws.Range("L132").Value = 11438.25
$ws.Range("M132").Value = -18858.7139
$ws.Range("N132").Value = -16498.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1420.8
$ws.Range("I68").Value = 1302
$ws.Range("K68").Value = 1302
$ws.Range("M68").Value = -553
$ws.Range("H71").Value = 1420.8
$ws.Range("I71").Value = 1302
$ws.Range("K71").Value = 6510
$ws.Range("M71").Value = -2766
$ws.Range("H82").Value = 2266.0667
$ws.Range("I82").Value = 2470.7144
$ws.Range("J82").Value = 2087
$ws.Range("K82").Value = 2470.7144
$ws.Range("L82").Value = 2087
$ws.Range("M82").Value = -2109.7144
$ws.Range("N82").Value = -2809
$ws.Range("H85").Value = 2266.0667
$ws.Range("I85").Value = 2470.7144
$ws.Range("J85").Value = 2087
$ws.Range("K85").Value = 2470.7144
$ws.Range("L85").Value = 2087
$ws.Range("M85").Value = -1222.7144
